$d = $word.ActiveDocument

$pairs = @(
    ,@("17+78=95", "16+57=73")
    ,@("1+79=80", "45+30=75")
    ,@("35+3=38", "25+68=93")
    ,@("81-68=13", "43+55=98")
    ,@("97-67=30", "30+26=56")
    ,@("99-35=64", "0+38=38")
    ,@("78-16=62", "55-54=1")
    ,@("29+20=49", "36+50=86")
    ,@("12+44=56", "69+11=80")
    ,@("84+13=97", "69-46=23")
    ,@("91-83=8", "2+29=31")
    ,@("45-13=32", "17-3=14")
    ,@("10+48=58", "36-6=30")
    ,@("33+19=52", "47-33=14")
    ,@("57-5=52", "19+41=60")
    ,@("27+0=27", "49+2=51")
    ,@("90-80=10", "77-36=41")
    ,@("3-1=2", "10-0=10")
    ,@("51-7=44", "21+8=29")
    ,@("51+47=98", "52-0=52")
    ,@("37+44=81", "41+7=48")
    ,@("48+15=63", "18+10=28")
    ,@("9+45=54", "43-36=7")
    ,@("49-6=43", "29-8=21")
    ,@("9+47=56", "27+51=78")
    ,@("7+7=14", "24-15=9")
    ,@("48+10=58", "64+24=88")
    ,@("42+53=95", "43-21=22")
    ,@("11+6=17", "91+4=95")
    ,@("38-23=15", "63-24=39")
    ,@("44-39=5", "4+50=54")
    ,@("58-31=27", "22+60=82")
    ,@("37+14=51", "30-15=15")
    ,@("56-41=15", "75+22=97")
    ,@("10+0=10", "45-39=6")
    ,@("64-3=61", "36-34=2")
    ,@("83-13=70", "76-25=51")
    ,@("65-2=63", "21+65=86")
    ,@("93-82=11", "26-1=25")
    ,@("14+21=35", "69-35=34")
    ,@("23+17=40", "86-83=3")
    ,@("8+84=92", "72-12=60")
    ,@("0+78=78", "52+37=89")
    ,@("86-35=51", "34+27=61")
    ,@("34+42=76", "47+16=63")
    ,@("18+76=94", "84-47=37")
    ,@("64-47=17", "52-23=29")
    ,@("67+16=83", "19+37=56")
    ,@("47-11=36", "72-48=24")
    ,@("15+19=34", "63+28=91")
    ,@("58-13=45", "74+19=93")
    ,@("54+21=75", "81-58=23")
    ,@("39-31=8", "30+21=51")
    ,@("5+69=74", "96-64=32")
    ,@("31-11=20", "78-37=41")
    ,@("66-2=64", "4+2=6")
    ,@("93-60=33", "21+77=98")
    ,@("97-48=49", "34+14=48")
    ,@("15+10=25", "14+70=84")
    ,@("53+19=72", "88-42=46")
    ,@("89-46=43", "8+38=46")
    ,@("57+17=74", "41+10=51")
    ,@("16+62=78", "55+2=57")
    ,@("2+97=99", "64+15=79")
    ,@("55+32=87", "59+28=87")
    ,@("36+17=53", "59+31=90")
    ,@("67+8=75", "47-12=35")
    ,@("22+71=93", "90+9=99")
    ,@("45+5=50", "90-70=20")
    ,@("80-36=44", "50-36=14")
    ,@("20+21=41", "12-5=7")
    ,@("35+12=47", "41+12=53")
    ,@("79-36=43", "56-53=3")
    ,@("22+15=37", "17+40=57")
    ,@("13+50=63", "75-25=50")
    ,@("18+28=46", "7+10=17")
    ,@("14+49=63", "30-6=24")
    ,@("51-47=4", "14+17=31")
    ,@("26+67=93", "58+34=92")
    ,@("83-64=19", "96-8=88")
    ,@("40+57=97", "37+18=55")
    ,@("67+14=81", "21+52=73")
    ,@("83-32=51", "12+51=63")
    ,@("97-15=82", "31+28=59")
    ,@("64-45=19", "83-22=61")
    ,@("49+50=99", "39-6=33")
    ,@("17+67=84", "44+33=77")
    ,@("14-1=13", "84-50=34")
    ,@("70+25=95", "69+2=71")
    ,@("47-37=10", "13+51=64")
    ,@("6+29=35", "27-2=25")
    ,@("23+58=81", "27-0=27")
    ,@("43-13=30", "96-4=92")
    ,@("35+9=44", "28-12=16")
    ,@("96-85=11", "15+40=55")
    ,@("76-31=45", "13+30=43")
    ,@("72-61=11", "0+34=34")
    ,@("17+2=19", "69+2=71")
    ,@("25-18=7", "83-59=24")
    ,@("48-23=25", "85-37=48")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

Write-Host "Done"
